$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 744.125
$ws.Range("I41").Value = 100.5
$ws.Range("J41").Value = 958.6667
$ws.Range("K41").Value = 100.5
$ws.Range("L41").Value = 958.6667
$ws.Range("M41").Value = 339.5
$ws.Range("N41").Value = -1838.6667
# Row 96
$ws.Range("H96").Value = 55562004
$ws.Range("I96").Value = 3308.2727
$ws.Range("J96").Value = 142868530
$ws.Range("K96").Value = 9924.8181
$ws.Range("L96").Value = 428605590
$ws.Range("M96").Value = -8551.8181
$ws.Range("N96").Value = -428608336
# Row 100
$ws.Range("H100").Value = 1566.9445
$ws.Range("I100").Value = 1586.0714
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1586.0714
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1045.0714
$ws.Range("N100").Value = -2582

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2213.9
$ws.Range("I2").Value = 2357.1428
$ws.Range("J2").Value = 1879.6666
$ws.Range("K2").Value = 2357.1428
$ws.Range("L2").Value = 1879.6666
$ws.Range("M2").Value = -2244.1428
$ws.Range("N2").Value = -2105.6666
# Row 32
$ws.Range("H32").Value = 35267.98
$ws.Range("I32").Value = 36108.266
$ws.Range("J32").Value = 27705.4
$ws.Range("K32").Value = 36108.266
$ws.Range("L32").Value = 27705.4
$ws.Range("M32").Value = -35821.266
$ws.Range("N32").Value = -28279.4
# Row 42
$ws.Range("H42").Value = 18996.666
$ws.Range("J42").Value = 18996.666
$ws.Range("L42").Value = 18996.666
$ws.Range("N42").Value = -19968.666
# Row 59
$ws.Range("H59").Value = 41000
$ws.Range("J59").Value = 41000
$ws.Range("L59").Value = 41000
$ws.Range("N59").Value = -42608
# Row 61
$ws.Range("H61").Value = 3009.6956
$ws.Range("I61").Value = 1765.3043
$ws.Range("J61").Value = 4254.087
$ws.Range("K61").Value = 1765.3043
$ws.Range("L61").Value = 4254.087
$ws.Range("M61").Value = -1553.3043
$ws.Range("N61").Value = -4678.087
# Row 92
$ws.Range("H92").Value = 23650
$ws.Range("J92").Value = 23650
$ws.Range("L92").Value = 23650
$ws.Range("N92").Value = -28642
# Row 102
$ws.Range("H102").Value = 55568950
$ws.Range("I102").Value = 62502256
$ws.Range("J102").Value = 102515.5
$ws.Range("K102").Value = 62502256
$ws.Range("L102").Value = 102515.5
$ws.Range("M102").Value = -62500634
$ws.Range("N102").Value = -105759.5
# Row 116
$ws.Range("H116").Value = 2213.9
$ws.Range("I116").Value = 2357.1428
$ws.Range("J116").Value = 1879.6666
$ws.Range("K116").Value = 2357.1428
$ws.Range("L116").Value = 1879.6666
$ws.Range("M116").Value = -63.14280000000008
$ws.Range("N116").Value = -6467.6666
# Row 122
$ws.Range("H122").Value = 2079.95
$ws.Range("I122").Value = 2287.92
$ws.Range("J122").Value = 1733.3334
$ws.Range("K122").Value = 6863.76
$ws.Range("L122").Value = 5200.0002
$ws.Range("M122").Value = -4413.76
$ws.Range("N122").Value = -10100.0002
# Row 136
$ws.Range("H136").Value = 3009.6956
$ws.Range("I136").Value = 1765.3043
$ws.Range("J136").Value = 4254.087
$ws.Range("K136").Value = 5295.9129
$ws.Range("L136").Value = 12762.261
$ws.Range("M136").Value = -2745.9129
$ws.Range("N136").Value = -17862.261

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2213.9
$ws.Range("I3").Value = 2357.1428
$ws.Range("J3").Value = 1879.6666
$ws.Range("K3").Value = 2357.1428
$ws.Range("L3").Value = 1879.6666
$ws.Range("M3").Value = -2243.1428
$ws.Range("N3").Value = -2107.6666
# Row 92
$ws.Range("H92").Value = 45464.332
$ws.Range("J92").Value = 45464.332
$ws.Range("L92").Value = 45464.332
$ws.Range("N92").Value = -50456.332

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 2316
$ws.Range("I99").Value = 2253
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2253
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -755
$ws.Range("N99").Value = -5396
# Row 103
$ws.Range("H103").Value = 17908.334
$ws.Range("J103").Value = 19990
$ws.Range("L103").Value = 19990
$ws.Range("N103").Value = -22334
# Row 122
$ws.Range("H122").Value = 75715.56
$ws.Range("I122").Value = 86432.07000000001
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 259296.21
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -256846.21
$ws.Range("N122").Value = -7000
# Row 126
$ws.Range("H126").Value = 2316
$ws.Range("I126").Value = 2253
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 6759
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -4289
$ws.Range("N126").Value = -12140

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 2062
$ws.Range("I63").Value = 110
$ws.Range("J63").Value = 4014
$ws.Range("K63").Value = 330
$ws.Range("L63").Value = 12042
$ws.Range("M63").Value = 419
$ws.Range("N63").Value = -13540
# Row 66
$ws.Range("H66").Value = 2062
$ws.Range("I66").Value = 110
$ws.Range("J66").Value = 4014
$ws.Range("K66").Value = 990
$ws.Range("L66").Value = 36126
$ws.Range("M66").Value = 2754
$ws.Range("N66").Value = -43614
# Row 107
$ws.Range("H107").Value = 40501.2
$ws.Range("I107").Value = 33834.332
$ws.Range("J107").Value = 50501.5
$ws.Range("K107").Value = 101502.996
$ws.Range("L107").Value = 151504.5
$ws.Range("M107").Value = -99582.99600000001
$ws.Range("N107").Value = -155344.5
# Row 131
$ws.Range("H131").Value = 2213.1396
$ws.Range("I131").Value = 17153.166
$ws.Range("J131").Value = 1092.6375
$ws.Range("K131").Value = 51459.49800000001
$ws.Range("L131").Value = 3277.9125
$ws.Range("M131").Value = -46419.49800000001
$ws.Range("N131").Value = -13357.9125

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 26000
$ws.Range("J58").Value = 26000
$ws.Range("L58").Value = 26000
$ws.Range("N58").Value = -26554
# Row 80
$ws.Range("H80").Value = 458018.62
$ws.Range("J80").Value = 4625
$ws.Range("L80").Value = 4625
$ws.Range("N80").Value = -6621
# Row 83
$ws.Range("H83").Value = 458018.62
$ws.Range("J83").Value = 4625
$ws.Range("L83").Value = 23125
$ws.Range("N83").Value = -33109
# Row 102
$ws.Range("H102").Value = 2394.7693
$ws.Range("I102").Value = 1979.7778
$ws.Range("J102").Value = 3328.5
$ws.Range("K102").Value = 1979.7778
$ws.Range("L102").Value = 3328.5
$ws.Range("M102").Value = -357.7778000000001
$ws.Range("N102").Value = -6572.5
# Row 130
$ws.Range("H130").Value = 47189.6
$ws.Range("J130").Value = 47189.6
$ws.Range("L130").Value = 47189.6
$ws.Range("N130").Value = -57229.6

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 26664
$ws.Range("J6").Value = 29996
$ws.Range("L6").Value = 29996
$ws.Range("N6").Value = -30220
# Row 57
$ws.Range("H57").Value = 77777
$ws.Range("J57").Value = 77777
$ws.Range("L57").Value = 77777
$ws.Range("N57").Value = -78909
# Row 100
$ws.Range("H100").Value = 2252.2
$ws.Range("I100").Value = 1946.8889
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1946.8889
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1405.8889
$ws.Range("N100").Value = -6082

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 27626.025
$ws.Range("I136").Value = 100947.5
$ws.Range("J136").Value = 2342.7585
$ws.Range("K136").Value = 302842.5
$ws.Range("L136").Value = 7028.2755
$ws.Range("M136").Value = -300292.5
$ws.Range("N136").Value = -12128.2755
